# Scheduled-runner style update of Leve profit-tracking sheets (ALC/BSM/CRP/CUL/LTW/WVR).
# Refreshes cached market-board derived figures (average prices, leve costs, profits)
# for a set of specific leve rows per sheet, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 232.1
$ws.Range("I103").Value = 204
$ws.Range("J103").Value = 266.44446
$ws.Range("K103").Value = 612
$ws.Range("L103").Value = 799.33338
$ws.Range("M103").Value = -26
$ws.Range("N103").Value = -1971.33338
$ws.Range("H116").Value = 2312.1765
$ws.Range("I116").Value = 1817.5
$ws.Range("J116").Value = 2582
$ws.Range("K116").Value = 1817.5
$ws.Range("L116").Value = 2582
$ws.Range("M116").Value = 1624.5
$ws.Range("N116").Value = -9466
$ws.Range("H138").Value = 5885350.5
$ws.Range("I138").Value = 6280
$ws.Range("J138").Value = 6898983.5
$ws.Range("K138").Value = 18840
$ws.Range("L138").Value = 20696950.5
$ws.Range("M138").Value = -13700
$ws.Range("N138").Value = -20707230.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 22478.416
$ws.Range("I82").Value = 3757
$ws.Range("J82").Value = 31839.125
$ws.Range("K82").Value = 3757
$ws.Range("L82").Value = 31839.125
$ws.Range("M82").Value = -3374
$ws.Range("N82").Value = -32605.125
$ws.Range("H85").Value = 22478.416
$ws.Range("I85").Value = 3757
$ws.Range("J85").Value = 31839.125
$ws.Range("K85").Value = 3757
$ws.Range("L85").Value = 31839.125
$ws.Range("M85").Value = -2431
$ws.Range("N85").Value = -34491.125
$ws.Range("H105").Value = 1793
$ws.Range("I105").Value = 1638.0769
$ws.Range("K105").Value = 1638.0769
$ws.Range("M105").Value = 108.9231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5558273
$ws.Range("I31").Value = 2845.85
$ws.Range("J31").Value = 12502557
$ws.Range("K31").Value = 2845.85
$ws.Range("L31").Value = 12502557
$ws.Range("M31").Value = -2550.85
$ws.Range("N31").Value = -12503147
$ws.Range("H34").Value = 5558273
$ws.Range("I34").Value = 2845.85
$ws.Range("J34").Value = 12502557
$ws.Range("K34").Value = 2845.85
$ws.Range("L34").Value = 12502557
$ws.Range("M34").Value = -2643.85
$ws.Range("N34").Value = -12502961
$ws.Range("H132").Value = 2287
$ws.Range("I132").Value = 1968.2142
$ws.Range("J132").Value = 3774.6667
$ws.Range("K132").Value = 5904.642599999999
$ws.Range("L132").Value = 11324.0001
$ws.Range("M132").Value = -3374.642599999999
$ws.Range("N132").Value = -16384.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 356.2
$ws.Range("I7").Value = 245.25
$ws.Range("J7").Value = 800
$ws.Range("K7").Value = 735.75
$ws.Range("L7").Value = 2400
$ws.Range("M7").Value = -623.75
$ws.Range("N7").Value = -2624
$ws.Range("H33").Value = 7668.5
$ws.Range("I33").Value = 14958.333
$ws.Range("J33").Value = 378.66666
$ws.Range("K33").Value = 89749.99800000001
$ws.Range("L33").Value = 2271.99996
$ws.Range("M33").Value = -89466.99800000001
$ws.Range("N33").Value = -2837.99996
$ws.Range("H47").Value = 400.5
$ws.Range("I47").Value = 280.6
$ws.Range("J47").Value = 1000
$ws.Range("K47").Value = 841.8000000000001
$ws.Range("L47").Value = 3000
$ws.Range("M47").Value = -410.8000000000001
$ws.Range("N47").Value = -3862
$ws.Range("H49").Value = 3625
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 3625
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 10875
$ws.Range("M49").ClearContents()
$ws.Range("N49").Value = -11187
$ws.Range("H68").Value = 2161.6736
$ws.Range("I68").Value = 1560.5209
$ws.Range("J68").Value = 2738.78
$ws.Range("K68").Value = 4681.5627
$ws.Range("L68").Value = 8216.34
$ws.Range("M68").Value = -3870.5627
$ws.Range("N68").Value = -9838.34
$ws.Range("H71").Value = 2161.6736
$ws.Range("I71").Value = 1560.5209
$ws.Range("J71").Value = 2738.78
$ws.Range("K71").Value = 14044.6881
$ws.Range("L71").Value = 24649.02
$ws.Range("M71").Value = -9988.688099999999
$ws.Range("N71").Value = -32761.02
$ws.Range("H131").Value = 3835823.5
$ws.Range("J131").Value = 7937460.5
$ws.Range("L131").Value = 23812381.5
$ws.Range("N131").Value = -23822461.5
$ws.Range("H132").Value = 1219.2307
$ws.Range("I132").Value = 453.84616
$ws.Range("K132").Value = 4084.61544
$ws.Range("M132").Value = -1554.61544

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5610.5283
$ws.Range("I132").Value = 8368.585999999999
$ws.Range("J132").Value = 2277.875
$ws.Range("K132").Value = 25105.758
$ws.Range("L132").Value = 6833.625
$ws.Range("M132").Value = -22575.758
$ws.Range("N132").Value = -11893.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 54684.832
$ws.Range("J111").Value = 54684.832
$ws.Range("L111").Value = 54684.832
$ws.Range("N111").Value = -62864.832
$ws.Range("H119").Value = 32155.9
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 32155.9
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 32155.9
$ws.Range("N119").Value = -41831.9
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("H121").Value = 19000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 19000
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 19000
$ws.Range("N121").Value = -22494
$ws.Range("H122").Value = 1149.909
$ws.Range("I122").Value = 1117.6666
$ws.Range("J122").Value = 1295
$ws.Range("K122").Value = 3352.9998
$ws.Range("L122").Value = 3885
$ws.Range("M122").Value = -902.9998000000001
$ws.Range("N122").Value = -8785
$ws.Range("H123").Value = 52194
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 52194
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 52194
$ws.Range("N123").Value = -61994
$ws.Range("H124").Value = 134714.5
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 134714.5
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 134714.5
$ws.Range("N124").Value = -144534.5
$ws.Range("H125").Value = 52666.668
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 52666.668
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 52666.668
$ws.Range("N125").Value = -62506.668
$ws.Range("H126").Value = 794.4211
$ws.Range("I126").Value = 661.8461
$ws.Range("J126").Value = 1081.6666
$ws.Range("K126").Value = 1985.5383
$ws.Range("L126").Value = 3244.9998
$ws.Range("M126").Value = 484.4617000000001
$ws.Range("N126").Value = -8184.9998
$ws.Range("H127").Value = 0
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H128").Value = 63805
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 63805
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 63805
$ws.Range("N128").Value = -73765
$ws.Range("H129").Value = 50000
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 50000
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 50000
$ws.Range("N129").Value = -60000
$ws.Range("H130").Value = 41500
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 41500
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 41500
$ws.Range("N130").Value = -51540
$ws.Range("H131").Value = 46500
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 46500
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 46500
$ws.Range("N131").Value = -56580
$ws.Range("H132").Value = 1131.5636
$ws.Range("I132").Value = 903.3095
$ws.Range("J132").Value = 1869
$ws.Range("K132").Value = 2709.9285
$ws.Range("L132").Value = 5607
$ws.Range("M132").Value = -179.9285
$ws.Range("N132").Value = -10667
$ws.Range("H133").Value = 44667.832
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 44667.832
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 44667.832
$ws.Range("N133").Value = -54787.832
$ws.Range("H135").Value = 53591.8
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 53591.8
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 53591.8
$ws.Range("N135").Value = -63731.8
$ws.Range("H136").Value = 5901.48
$ws.Range("I136").Value = 8074
$ws.Range("J136").Value = 1284.875
$ws.Range("K136").Value = 24222
$ws.Range("L136").Value = 3854.625
$ws.Range("M136").Value = -21672
$ws.Range("N136").Value = -8954.625
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 79209.336
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 79209.336
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 79209.336
$ws.Range("N138").Value = -89489.336
$ws.Range("H139").Value = 46759
$ws.Range("I139").Value = 40650
$ws.Range("J139").Value = 48286.25
$ws.Range("K139").Value = 40650
$ws.Range("L139").Value = 48286.25
$ws.Range("M139").Value = -35510
$ws.Range("N139").Value = -58566.25
$ws.Range("H140").Value = 57929
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 57929
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 57929
$ws.Range("N140").Value = -68289
$ws.Range("H141").Value = 68715
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 68715
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 68715
$ws.Range("N141").Value = -79075
